$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "17.1.2 Доля национального бюджета, финансируемая за счет внутренних налогов"
$ws.Range("B4").Select()
